# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled
# update). Price (column D) and Volume(1h) (column E) values are updated
# for most rows; rows 32/33 additionally swap the Stellar / ImmutableX
# entries (name, link, price, volume) as the ranking order changed.
#
# Column D values are prefixed with a leading apostrophe so Excel stores
# them as literal text (matching the workbook's original inline-string
# representation) instead of re-parsing them as numbers, which would
# silently drop meaningful trailing zeros (e.g. "111.00" -> 111,
# "0.06920" -> 0.0692).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.132.31"
$ws.Range("E2").Value = "  -2.95%  "

$ws.Range("D3").Value = "'1.911.52"
$ws.Range("E3").Value = "  -3.70%  "

$ws.Range("E4").Value = "  -1.36%  "

$ws.Range("D5").Value = "'328.17"
$ws.Range("E5").Value = "  -0.47%  "

$ws.Range("E6").Value = "  -1.24%  "

$ws.Range("D7").Value = "'0.4640"
$ws.Range("E7").Value = "  -5.65%  "

$ws.Range("D8").Value = "'0.4007"
$ws.Range("E8").Value = "  -3.59%  "

$ws.Range("D9").Value = "'53.25"
$ws.Range("E9").Value = "  -3.23%  "

$ws.Range("E10").Value = "  -4.83%  "

$ws.Range("D11").Value = "'1.042"
$ws.Range("E11").Value = "  -3.65%  "

$ws.Range("D12").Value = "'21.95"
$ws.Range("E12").Value = "  -3.23%  "

$ws.Range("D13").Value = "'1.903.15"
$ws.Range("E13").Value = "  -8.12%  "

$ws.Range("D14").Value = "'7.413"
$ws.Range("E14").Value = "  -6.20%  "

$ws.Range("D15").Value = "'6.055"
$ws.Range("E15").Value = "  -4.84%  "

$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -1.39%  "

$ws.Range("D17").Value = "'89.56"
$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("E18").Value = "  -2.91%  "

$ws.Range("D19").Value = "'0.06604"
$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").Value = "'17.86"
$ws.Range("E20").Value = "  -7.17%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  -1.14%  "

$ws.Range("D22").Value = "'5.745"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").Value = "'28.116.46"
$ws.Range("E23").Value = "  -3.19%  "

$ws.Range("D24").Value = "'11.17"
$ws.Range("E24").Value = "  -5.91%  "

$ws.Range("D25").Value = "'2.305"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "'2.113.11"
$ws.Range("E26").Value = "  -7.56%  "

$ws.Range("D27").Value = "'152.86"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("E28").Value = "  -3.08%  "

$ws.Range("D29").Value = "'5.768"
$ws.Range("E29").Value = "  -6.78%  "

$ws.Range("D30").Value = "'2.127"
$ws.Range("E30").Value = "  -4.07%  "

$ws.Range("D31").Value = "'123.58"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32/33: Stellar and ImmutableX swap places in the ranking.
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09644"
$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.9737"
$ws.Range("E33").Value = "  -5.77%  "

$ws.Range("D34").Value = "'1.449"
$ws.Range("E34").Value = "  -4.11%  "

$ws.Range("D35").Value = "'5.552"
$ws.Range("E35").Value = "  -4.26%  "

$ws.Range("D36").Value = "'3.633"
$ws.Range("E36").Value = "  -2.81%  "

$ws.Range("D37").Value = "'1.277"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").Value = "'8.785"
$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("D39").Value = "'0.02293"
$ws.Range("E39").Value = "  -4.31%  "

$ws.Range("D40").Value = "'0.06136"
$ws.Range("E40").Value = "  -3.41%  "

$ws.Range("D41").Value = "'0.6154"
$ws.Range("E41").Value = "  -4.32%  "

$ws.Range("D42").Value = "'10.92"
$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -1.29%  "

$ws.Range("D44").Value = "'0.1906"
$ws.Range("E44").Value = "  -3.19%  "

$ws.Range("E45").Value = "  -4.10%  "

$ws.Range("D46").Value = "'0.5859"
$ws.Range("E46").Value = "  -4.60%  "

$ws.Range("D47").Value = "'12.77"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "'2.015"
$ws.Range("E48").Value = "  -6.00%  "

$ws.Range("D49").Value = "'3.433"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").Value = "'0.06920"
$ws.Range("E50").Value = "  -0.02%  "

$ws.Range("D51").Value = "'111.00"
$ws.Range("E51").Value = "  -1.42%  "
